$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fill: new columns H, J:Q (additional run data) and R (labels a/b) ---
$ws.Range("H3").Value = 22.541966326743701
$ws.Range("J3").Value = 22.429426987332899
$ws.Range("K3").Value = 22.3846201831707
$ws.Range("L3").Value = 22.638872942387099
$ws.Range("M3").Value = 22.429426987332899
$ws.Range("N3").Value = 22.641297160130001
$ws.Range("O3").Value = 22.429426987332899
$ws.Range("P3").Value = 22.3846201831707
$ws.Range("Q3").Value = 22.4283245797067
$ws.Range("H4").Value = 9.7972494149063394
$ws.Range("J4").Value = 9.7972494149063394
$ws.Range("K4").Value = 9.7972494149063394
$ws.Range("L4").Value = 9.7972494149063394
$ws.Range("M4").Value = 9.7981551468909593
$ws.Range("N4").Value = 9.7972494149063394
$ws.Range("O4").Value = 9.7972494149063394
$ws.Range("P4").Value = 9.7972494149063394
$ws.Range("Q4").Value = 9.7972494149063394
$ws.Range("R4").Value = "a"
$ws.Range("H5").Value = 1.12140865001139
$ws.Range("J5").Value = 1.12140865001139
$ws.Range("K5").Value = 1.12140865001139
$ws.Range("L5").Value = 1.12140865001139
$ws.Range("M5").Value = 1.12140865001139
$ws.Range("N5").Value = 1.12140865001139
$ws.Range("O5").Value = 1.12140865001139
$ws.Range("P5").Value = 1.12140865001139
$ws.Range("Q5").Value = 1.12140865001139
$ws.Range("R5").Value = "a"
$ws.Range("H6").Value = 0.33638761656947502
$ws.Range("J6").Value = 0.33638761656947502
$ws.Range("K6").Value = 0.345569198312149
$ws.Range("L6").Value = 0.345569198312149
$ws.Range("M6").Value = 0.345569198312149
$ws.Range("N6").Value = 0.33638761656947502
$ws.Range("O6").Value = 0.36904471329701399
$ws.Range("P6").Value = 0.345569198312149
$ws.Range("Q6").Value = 0.33704147201486401
$ws.Range("H7").Value = 0.55971682890649099
$ws.Range("J7").Value = 0.55971682890649099
$ws.Range("K7").Value = 0.55971682890649099
$ws.Range("L7").Value = 0.55971682890649099
$ws.Range("M7").Value = 0.55971682890649099
$ws.Range("N7").Value = 0.55971682890649099
$ws.Range("O7").Value = 0.55971682890649099
$ws.Range("P7").Value = 0.55971682890649099
$ws.Range("Q7").Value = 0.55971682890649099
$ws.Range("R7").Value = "a"
$ws.Range("H8").Value = 6.69496542660006
$ws.Range("J8").Value = 6.8417570676701596
$ws.Range("K8").Value = 6.69496542660006
$ws.Range("L8").Value = 6.3843814412801896
$ws.Range("M8").Value = 6.8417570676701596
$ws.Range("N8").Value = 6.69496542660006
$ws.Range("O8").Value = 6.8417570676701596
$ws.Range("P8").Value = 6.8417570676701596
$ws.Range("Q8").Value = 6.3843814412801896
$ws.Range("R8").Value = "b"
$ws.Range("H9").Value = 9.0426940999840308
$ws.Range("J9").Value = 9.0426940999840308
$ws.Range("K9").Value = 9.7293938006062994
$ws.Range("L9").Value = 9.7293938006062994
$ws.Range("M9").Value = 18.2976766853043
$ws.Range("N9").Value = 9.0426940999840308
$ws.Range("O9").Value = 9.7293938006062994
$ws.Range("P9").Value = 9.7293938006062994
$ws.Range("Q9").Value = 9.0426940999840308
$ws.Range("R9").Value = "b"
$ws.Range("H10").Value = 8.70883502813383
$ws.Range("J10").Value = 9.0694649952438393
$ws.Range("K10").Value = 10.577553829768499
$ws.Range("L10").Value = 8.70883502813383
$ws.Range("M10").Value = 10.578650337052
$ws.Range("N10").Value = 8.7088092889942992
$ws.Range("O10").Value = 8.70883502813383
$ws.Range("P10").Value = 9.1103212675399501
$ws.Range("Q10").Value = 9.0694649952438393
$ws.Range("H11").Value = 6.8998146794585802
$ws.Range("J11").Value = 7.3319909444984699
$ws.Range("K11").Value = 6.8479746403691104
$ws.Range("L11").Value = 6.8998146794585802
$ws.Range("M11").Value = 7.3319909444984699
$ws.Range("N11").Value = 6.8998146794585802
$ws.Range("O11").Value = 6.8998146794585802
$ws.Range("P11").Value = 6.8479746403691104
$ws.Range("Q11").Value = 7.3319909444984699
$ws.Range("R11").Value = "b"
$ws.Range("H12").Value = 6.0671907663239697
$ws.Range("J12").Value = 6.3802043727122602
$ws.Range("K12").Value = 6.3802043727122602
$ws.Range("L12").Value = 6.3735141686971799
$ws.Range("M12").Value = 6.0671907663239697
$ws.Range("N12").Value = 6.3802043727122602
$ws.Range("O12").Value = 10.5265348463649
$ws.Range("P12").Value = 6.0671907663239697
$ws.Range("Q12").Value = 6.8869665858309697
$ws.Range("H13").Value = 6.1905284490298103
$ws.Range("J13").Value = 6.1905284490298103
$ws.Range("K13").Value = 6.1905284490298103
$ws.Range("L13").Value = 6.7161736385278097
$ws.Range("M13").Value = 6.1766118569259802
$ws.Range("N13").Value = 6.1766118569259802
$ws.Range("O13").Value = 6.1905284490298103
$ws.Range("P13").Value = 6.1905284490298103
$ws.Range("Q13").Value = 6.1808571887563399
$ws.Range("H14").Value = 7.6038071444614399
$ws.Range("J14").Value = 7.5474563295486101
$ws.Range("K14").Value = 7.6761816308303796
$ws.Range("L14").Value = 7.5652344489831904
$ws.Range("M14").Value = 8.0025282382155503
$ws.Range("N14").Value = 7.5379780855500904
$ws.Range("O14").Value = 7.52074384239359
$ws.Range("P14").Value = 7.5597768681289299
$ws.Range("Q14").Value = 7.6038071444614399
$ws.Range("H15").Value = 0.29195881139122198
$ws.Range("J15").Value = 0.29052642476740298
$ws.Range("K15").Value = 0.28947602411886503
$ws.Range("L15").Value = 0.29195881139122198
$ws.Range("M15").Value = 0.29198208906158901
$ws.Range("N15").Value = 0.29398060498425099
$ws.Range("O15").Value = 0.30230228130627501
$ws.Range("P15").Value = 0.31508403032235799
$ws.Range("Q15").Value = 0.29402813796384702
$ws.Range("H16").Value = 0.54217027264039896
$ws.Range("J16").Value = 0.54217027264039896
$ws.Range("K16").Value = 0.54217027264039896
$ws.Range("L16").Value = 0.54217027264039896
$ws.Range("M16").Value = 0.54217027264039896
$ws.Range("N16").Value = 0.54217027264039896
$ws.Range("O16").Value = 0.54217027264039896
$ws.Range("P16").Value = 0.54217027264039896
$ws.Range("Q16").Value = 0.54217027264039896
$ws.Range("R16").Value = "a"
$ws.Range("H17").Value = 0.73854470165158903
$ws.Range("J17").Value = 0.73854470165158903
$ws.Range("K17").Value = 0.73854470165158903
$ws.Range("L17").Value = 0.73483909231449596
$ws.Range("M17").Value = 0.73483909231449596
$ws.Range("N17").Value = 0.73854470165158903
$ws.Range("O17").Value = 0.73854470165158903
$ws.Range("P17").Value = 0.73854470165158903
$ws.Range("Q17").Value = 0.73854470165158903
$ws.Range("H18").Value = 1.4832748222186301
$ws.Range("J18").Value = 2.1223069289018701
$ws.Range("K18").Value = 1.4832748222186301
$ws.Range("L18").Value = 1.4832748222186301
$ws.Range("M18").Value = 2.1223069289018701
$ws.Range("N18").Value = 2.1223069289018701
$ws.Range("O18").Value = 1.4832748222186301
$ws.Range("P18").Value = 1.4832748222186301
$ws.Range("Q18").Value = 1.4832748222186301
$ws.Range("H19").Value = 1.5652926334756601
$ws.Range("J19").Value = 1.56594670734242
$ws.Range("K19").Value = 1.5652926334756601
$ws.Range("L19").Value = 1.5652926334756601
$ws.Range("M19").Value = 1.56594670734242
$ws.Range("N19").Value = 1.5652926334756601
$ws.Range("O19").Value = 1.5652926334756601
$ws.Range("P19").Value = 1.5652926334756601
$ws.Range("Q19").Value = 1.56470230382949
$ws.Range("H20").Value = 3.1803719036714102
$ws.Range("J20").Value = 3.1803719036714102
$ws.Range("K20").Value = 4.03361824086823
$ws.Range("L20").Value = 3.3853840538239499
$ws.Range("M20").Value = 3.1803719036714102
$ws.Range("N20").Value = 3.3853840538239499
$ws.Range("O20").Value = 3.1803719036714102
$ws.Range("P20").Value = 3.3853840538239499
$ws.Range("Q20").Value = 3.1803719036714102
$ws.Range("H21").Value = 1.68893109303075
$ws.Range("J21").Value = 1.68893109303075
$ws.Range("K21").Value = 1.68893109303075
$ws.Range("L21").Value = 1.68893109303075
$ws.Range("M21").Value = 1.68893109303075
$ws.Range("N21").Value = 1.68893109303075
$ws.Range("O21").Value = 1.68893109303075
$ws.Range("P21").Value = 1.68893109303075
$ws.Range("Q21").Value = 1.68893109303075
$ws.Range("R21").Value = "a"
$ws.Range("H22").Value = 2.3116449243396202
$ws.Range("J22").Value = 2.3116449243396202
$ws.Range("K22").Value = 2.3116449243396202
$ws.Range("L22").Value = 2.3116449243396202
$ws.Range("M22").Value = 2.3116449243396202
$ws.Range("N22").Value = 2.3116449243396202
$ws.Range("O22").Value = 3.7296710262705699
$ws.Range("P22").Value = 2.3116449243396202
$ws.Range("Q22").Value = 2.3116449239999999
$ws.Range("Q22").Font.Color = 0
$ws.Range("H23").Value = 2.1208650192724199
$ws.Range("J23").Value = 2.1208650192724199
$ws.Range("K23").Value = 2.1208650192724199
$ws.Range("L23").Value = 2.1208650192724199
$ws.Range("M23").Value = 2.1208650192724199
$ws.Range("N23").Value = 2.1208650192724199
$ws.Range("O23").Value = 5.9180577874489098
$ws.Range("P23").Value = 2.1208650192724199
$ws.Range("Q23").Value = 2.1208650192724199
$ws.Range("H24").Value = 1.25831324316183
$ws.Range("J24").Value = 1.25831324316183
$ws.Range("K24").Value = 1.25831324316183
$ws.Range("L24").Value = 1.25831324316183
$ws.Range("M24").Value = 1.25831324316183
$ws.Range("N24").Value = 1.25831324316183
$ws.Range("O24").Value = 1.25831324316183
$ws.Range("P24").Value = 3.5987869723542998
$ws.Range("Q24").Value = 1.25831324316183
$ws.Range("R24").Value = "a"

# --- Sheet view: zoom + selection (matches the saved view state) ---
$excel.ActiveWindow.Zoom = 112
$ws.Range("R11").Select()
